$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3000
$ws.Range("I18").Value = 3000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = -2716
$ws.Range("H62").Value = 3667.3865
$ws.Range("I62").Value = 3051.8684
$ws.Range("K62").Value = 3051.8684
$ws.Range("M62").Value = -2427.8684
$ws.Range("H65").Value = 3667.3865
$ws.Range("I65").Value = 3051.8684
$ws.Range("K65").Value = 15259.342
$ws.Range("M65").Value = -12139.342
$ws.Range("H70").Value = 6636.909
$ws.Range("I70").Value = 6948.7
$ws.Range("J70").Value = 6377.0835
$ws.Range("K70").Value = 20846.1
$ws.Range("L70").Value = 19131.2505
$ws.Range("M70").Value = -20576.1
$ws.Range("N70").Value = -19671.2505
$ws.Range("H73").Value = 6636.909
$ws.Range("I73").Value = 6948.7
$ws.Range("J73").Value = 6377.0835
$ws.Range("K73").Value = 20846.1
$ws.Range("L73").Value = 19131.2505
$ws.Range("M73").Value = -19910.1
$ws.Range("N73").Value = -21003.2505
$ws.Range("H107").Value = 1652.76
$ws.Range("I107").Value = 1453.25
$ws.Range("K107").Value = 1453.25
$ws.Range("M107").Value = 466.75
$ws.Range("H111").Value = 612.2
$ws.Range("I111").Value = 591.3333
$ws.Range("K111").Value = 1773.9999
$ws.Range("M111").Value = 1293.0001
$ws.Range("H112").Value = 5169.303
$ws.Range("I112").Value = 1699.3334
$ws.Range("J112").Value = 5516.3
$ws.Range("K112").Value = 5098.0002
$ws.Range("L112").Value = 16548.9
$ws.Range("M112").Value = -3990.0002
$ws.Range("N112").Value = -18764.9
$ws.Range("H137").Value = 10786.765
$ws.Range("I137").Value = 13842.56
$ws.Range("J137").Value = 2298.4443
$ws.Range("K137").Value = 41527.68
$ws.Range("L137").Value = 6895.3329
$ws.Range("M137").Value = -38977.68
$ws.Range("N137").Value = -11995.3329
$ws.Range("H138").Value = 31758.205
$ws.Range("I138").Value = 2007.6923
$ws.Range("J138").Value = 128447.375
$ws.Range("K138").Value = 6023.0769
$ws.Range("L138").Value = 385342.125
$ws.Range("M138").Value = -883.0769
$ws.Range("N138").Value = -395622.125
$ws.Range("H140").Value = 97326.336
$ws.Range("J140").Value = 97326.336
$ws.Range("L140").Value = 97326.336
$ws.Range("N140").Value = -107686.336

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21467.117
$ws.Range("I32").Value = 23487.479
$ws.Range("K32").Value = 23487.479
$ws.Range("M32").Value = -23200.479
$ws.Range("H61").Value = 5753
$ws.Range("I61").Value = 897
$ws.Range("J61").Value = 10609
$ws.Range("K61").Value = 897
$ws.Range("L61").Value = 10609
$ws.Range("M61").Value = -685
$ws.Range("N61").Value = -11033
$ws.Range("H74").Value = 133593.98
$ws.Range("I74").Value = 162860.81
$ws.Range("J74").Value = 13274.777
$ws.Range("K74").Value = 162860.81
$ws.Range("L74").Value = 13274.777
$ws.Range("M74").Value = -161986.81
$ws.Range("N74").Value = -15022.777
$ws.Range("H77").Value = 133593.98
$ws.Range("I77").Value = 162860.81
$ws.Range("J77").Value = 13274.777
$ws.Range("K77").Value = 814304.05
$ws.Range("L77").Value = 66373.88499999999
$ws.Range("M77").Value = -809936.05
$ws.Range("N77").Value = -75109.88499999999
$ws.Range("H132").Value = 1232.4147
$ws.Range("I132").Value = 1069.5676
$ws.Range("K132").Value = 3208.7028
$ws.Range("M132").Value = -678.7028
$ws.Range("H136").Value = 5753
$ws.Range("I136").Value = 897
$ws.Range("J136").Value = 10609
$ws.Range("K136").Value = 2691
$ws.Range("L136").Value = 31827
$ws.Range("M136").Value = -141
$ws.Range("N136").Value = -36927

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4188.5557
$ws.Range("I105").Value = 3003.3333
$ws.Range("J105").Value = 4781.1665
$ws.Range("K105").Value = 3003.3333
$ws.Range("L105").Value = 4781.1665
$ws.Range("M105").Value = -1256.3333
$ws.Range("N105").Value = -8275.166499999999
$ws.Range("H107").Value = 872.1667
$ws.Range("I107").Value = 846.8
$ws.Range("K107").Value = 846.8
$ws.Range("M107").Value = 1073.2
$ws.Range("H134").Value = 2147.5667
$ws.Range("I134").Value = 1577.4231
$ws.Range("K134").Value = 4732.2693
$ws.Range("M134").Value = -2197.2693

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2009.909
$ws.Range("I16").Value = 1566.2222
$ws.Range("K16").Value = 1566.2222
$ws.Range("M16").Value = -1279.2222
$ws.Range("H31").Value = 3335582.2
$ws.Range("I31").Value = 4546568
$ws.Range("K31").Value = 4546568
$ws.Range("M31").Value = -4546273
$ws.Range("H34").Value = 3335582.2
$ws.Range("I34").Value = 4546568
$ws.Range("K34").Value = 4546568
$ws.Range("M34").Value = -4546366
$ws.Range("H86").Value = 13804.863
$ws.Range("I86").Value = 16490.416
$ws.Range("K86").Value = 16490.416
$ws.Range("M86").Value = -15367.416
$ws.Range("H89").Value = 13804.863
$ws.Range("I89").Value = 16490.416
$ws.Range("K89").Value = 82452.08
$ws.Range("M89").Value = -76836.08
$ws.Range("H94").Value = 1430.92
$ws.Range("I94").Value = 1122.5555
$ws.Range("J94").Value = 1604.375
$ws.Range("K94").Value = 1122.5555
$ws.Range("L94").Value = 1604.375
$ws.Range("M94").Value = -671.5554999999999
$ws.Range("N94").Value = -2506.375
$ws.Range("H113").Value = 2009.909
$ws.Range("I113").Value = 1566.2222
$ws.Range("K113").Value = 1566.2222
$ws.Range("M113").Value = 603.7778000000001
$ws.Range("H122").Value = 1452.8572
$ws.Range("I122").Value = 1462.88
$ws.Range("J122").Value = 1427.8
$ws.Range("K122").Value = 4388.64
$ws.Range("L122").Value = 4283.4
$ws.Range("M122").Value = -1938.64
$ws.Range("N122").Value = -9183.4
$ws.Range("H132").Value = 33430.773
$ws.Range("I132").Value = 40999.24
$ws.Range("J132").Value = 1895.5
$ws.Range("K132").Value = 122997.72
$ws.Range("L132").Value = 5686.5
$ws.Range("M132").Value = -120467.72
$ws.Range("N132").Value = -10746.5
$ws.Range("H134").Value = 2090.3333
$ws.Range("I134").Value = 1785.7693
$ws.Range("J134").Value = 2882.2
$ws.Range("K134").Value = 5357.3079
$ws.Range("L134").Value = 8646.599999999999
$ws.Range("M134").Value = -2822.3079
$ws.Range("N134").Value = -13716.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 559.8570999999999
$ws.Range("I98").Value = 298.25
$ws.Range("J98").Value = 664.5
$ws.Range("K98").Value = 894.75
$ws.Range("L98").Value = 1993.5
$ws.Range("M98").Value = 603.25
$ws.Range("N98").Value = -4989.5
$ws.Range("H107").Value = 1992.6818
$ws.Range("J107").Value = 793.6875
$ws.Range("L107").Value = 2381.0625
$ws.Range("N107").Value = -6221.0625
$ws.Range("H112").Value = 4243
$ws.Range("I112").Value = 3853.75
$ws.Range("K112").Value = 11561.25
$ws.Range("M112").Value = -10453.25
$ws.Range("H116").Value = 8889.9
$ws.Range("H133").Value = 4007.4546
$ws.Range("I133").Value = 3565.3
$ws.Range("K133").Value = 10695.9
$ws.Range("M133").Value = -5635.900000000001
$ws.Range("H139").Value = 2476.3333
$ws.Range("I139").Value = 2476.3333
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7428.999899999999
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2288.999899999999
$ws.Range("N139").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H113").Value = 2247.25
$ws.Range("I113").Value = 2170.8572
$ws.Range("J113").Value = 2354.2
$ws.Range("K113").Value = 2170.8572
$ws.Range("L113").Value = 2354.2
$ws.Range("M113").Value = -0.8571999999999207
$ws.Range("N113").Value = -6694.2
$ws.Range("H132").Value = 3333.9167
$ws.Range("I132").Value = 3000.8
$ws.Range("K132").Value = 9002.400000000001
$ws.Range("M132").Value = -6472.400000000001
$ws.Range("N26").ClearContents()
$ws.Range("N50").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2742
$ws.Range("I46").Value = 755.4286
$ws.Range("J46").Value = 5059.6665
$ws.Range("K46").Value = 755.4286
$ws.Range("L46").Value = 5059.6665
$ws.Range("M46").Value = -567.4286
$ws.Range("N46").Value = -5435.6665
$ws.Range("H68").Value = 3280.125
$ws.Range("I68").Value = 2776.4443
$ws.Range("K68").Value = 2776.4443
$ws.Range("M68").Value = -2027.4443
$ws.Range("H71").Value = 3280.125
$ws.Range("I71").Value = 2776.4443
$ws.Range("K71").Value = 13882.2215
$ws.Range("M71").Value = -10138.2215
$ws.Range("H119").Value = 101474.5
$ws.Range("J119").Value = 101474.5
$ws.Range("L119").Value = 101474.5
$ws.Range("N119").Value = -111150.5
$ws.Range("H136").Value = 3054.0789
$ws.Range("I136").Value = 2886.7
$ws.Range("K136").Value = 8660.099999999999
$ws.Range("M136").Value = -6110.099999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26357.477
$ws.Range("I136").Value = 27625.35
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 82876.04999999999
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -80326.04999999999
$ws.Range("N136").Value = -8100
